$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "-"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "[-, 'MCT-2A-CAD']"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[-, 'MCT-2A-CAD']"

# Row 6
$ws.Range("F6").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"

# Row 7
$ws.Range("C7").Value = "MCT-1A-Desenho Técnico"
$ws.Range("F7").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
